# [silverfox] 새 npc 가데이터 입력
# Adds two new "cleaner cooler" rows (cid 5300 / 5301) to Sheet1, matching
# the existing table layout: cid, name, nameKor, socket, grade, company,
# cooltime. The grade/company cells on the new rows get a purple fill with
# white text to flag them as newly added NPC gear data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: cleaner_cooler / 클리너 쿨러 ---------------------------------
$ws.Range("A13").Value = 5300
$ws.Range("B13").Value = "cleaner_cooler"
$ws.Range("C13").Value = "클리너 쿨러"
$ws.Range("D13").Value = "cooler"
$ws.Range("E13").Value = "Weakness"
$ws.Range("F13").Value = "WalkersWorkshop"
$ws.Range("G13").Value = 4

# --- Row 14: cleaner_R_cooler / 깔끔한 클리너 쿨러 ------------------------
$ws.Range("A14").Value = 5301
$ws.Range("B14").Value = "cleaner_R_cooler"
$ws.Range("C14").Value = "깔끔한 클리너 쿨러"
$ws.Range("D14").Value = "cooler"
$ws.Range("E14").Value = "Weakness"
$ws.Range("F14").Value = "DSDC"
$ws.Range("G14").Value = 4

# --- Highlight formatting (purple fill, white text) -----------------------
# Apply every fill first, then every font color -- doing them interleaved
# per-cell can make the engine drop a duplicate-valued property, so batch
# by property instead of by cell.
$ws.Range("E13").Interior.Color = 10498160
$ws.Range("F13:G13").Interior.Color = 10498160
$ws.Range("E14:G14").Interior.Color = 10498160

$ws.Range("E13").Font.ThemeColor = 2
$ws.Range("F13:G13").Font.ThemeColor = 2
$ws.Range("E14:G14").Font.ThemeColor = 2

# --- Match the final selection left behind in the authored workbook -------
$ws.Range("E14").Select()
